# Auto-generated edit script: applies the scheduled market-data refresh
# to the per-job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each row updates the market-price / profit columns (H-N) with refreshed
# Universalis price data; a few rows gain or lose a column entirely
# because the corresponding price was previously absent/zeroed out.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 2399.6667
$ws.Range("I40").Value = 1679.6
$ws.Range("K40").Value = 1679.6
$ws.Range("M40").Value = -1504.6
# Row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 5596.231
$ws.Range("I43").Value = 6464
$ws.Range("J43").Value = 4852.4287
$ws.Range("K43").Value = 6464
$ws.Range("L43").Value = 4852.4287
$ws.Range("M43").Value = -6395
$ws.Range("N43").Value = -4990.4287
# Row 76 (Leve Item ID 12602)
$ws.Range("H76").Value = 5203.2666
$ws.Range("I76").Value = 4736.7334
$ws.Range("J76").Value = 5669.8
$ws.Range("K76").Value = 4736.7334
$ws.Range("L76").Value = 5669.8
$ws.Range("M76").Value = -4421.7334
$ws.Range("N76").Value = -6299.8
# Row 79 (Leve Item ID 12602)
$ws.Range("H79").Value = 5203.2666
$ws.Range("I79").Value = 4736.7334
$ws.Range("J79").Value = 5669.8
$ws.Range("K79").Value = 4736.7334
$ws.Range("L79").Value = 5669.8
$ws.Range("M79").Value = -3644.7334
$ws.Range("N79").Value = -7853.8
# Row 96 (Leve Item ID 19894)
$ws.Range("H96").Value = 7143496.5
$ws.Range("I96").Value = 11905003
$ws.Range("J96").Value = 1237.25
$ws.Range("K96").Value = 35715009
$ws.Range("L96").Value = 3711.75
$ws.Range("M96").Value = -35713636
$ws.Range("N96").Value = -6457.75

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 14060.412
$ws.Range("I32").Value = 14267.322
$ws.Range("K32").Value = 14267.322
$ws.Range("M32").Value = -13980.322
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 4197.5
$ws.Range("I61").Value = 1996.1666
$ws.Range("K61").Value = 1996.1666
$ws.Range("M61").Value = -1784.1666
# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 1130.4286
$ws.Range("I110").Value = 1102
$ws.Range("K110").Value = 1102
$ws.Range("M110").Value = 943
# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 2199.3635
$ws.Range("I122").Value = 1919.3
$ws.Range("K122").Value = 5757.9
$ws.Range("M122").Value = -3307.9
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 4197.5
$ws.Range("I136").Value = 1996.1666
$ws.Range("K136").Value = 5988.4998
$ws.Range("M136").Value = -3438.4998

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (Leve Item ID 14149)
$ws.Range("H20").Value = 21334.16
$ws.Range("I20").Value = 32632.438
$ws.Range("J20").Value = 1248.3334
$ws.Range("K20").Value = 32632.438
$ws.Range("L20").Value = 1248.3334
$ws.Range("M20").Value = -32385.438
$ws.Range("N20").Value = -1742.3334
# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 2061.111
$ws.Range("I94").Value = 1380
$ws.Range("J94").Value = 3131.4285
$ws.Range("K94").Value = 1380
$ws.Range("L94").Value = 3131.4285
$ws.Range("M94").Value = -929
$ws.Range("N94").Value = -4033.4285
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 3856.5625
$ws.Range("I134").Value = 3856.5625
$ws.Range("K134").Value = 11569.6875
$ws.Range("M134").Value = -9034.6875

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 2324.25
$ws.Range("I16").Value = 1997
$ws.Range("K16").Value = 1997
$ws.Range("M16").Value = -1710
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2130423.2
$ws.Range("I31").Value = 2501936.5
$ws.Range("J31").Value = 7489.2856
$ws.Range("K31").Value = 2501936.5
$ws.Range("L31").Value = 7489.2856
$ws.Range("M31").Value = -2501641.5
$ws.Range("N31").Value = -8079.2856
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2130423.2
$ws.Range("I34").Value = 2501936.5
$ws.Range("J34").Value = 7489.2856
$ws.Range("K34").Value = 2501936.5
$ws.Range("L34").Value = 7489.2856
$ws.Range("M34").Value = -2501734.5
$ws.Range("N34").Value = -7893.2856
# Row 62 (Leve Item ID 12580)
$ws.Range("H62").Value = 9204.799999999999
$ws.Range("J62").Value = 10356
$ws.Range("L62").Value = 10356
$ws.Range("N62").Value = -11604
# Row 65 (Leve Item ID 12580)
$ws.Range("H65").Value = 9204.799999999999
$ws.Range("J65").Value = 10356
$ws.Range("L65").Value = 51780
$ws.Range("N65").Value = -58020
# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 2324.25
$ws.Range("I113").Value = 1997
$ws.Range("K113").Value = 1997
$ws.Range("M113").Value = 173

$ws = $wb.Worksheets.Item("CUL")
# Row 70 (Leve Item ID 12867)
$ws.Range("H70").Value = 4272.3335
$ws.Range("I70").Value = 2089.3333
$ws.Range("K70").Value = 6267.999899999999
$ws.Range("M70").Value = -5952.999899999999
# Row 73 (Leve Item ID 12867)
$ws.Range("H73").Value = 4272.3335
$ws.Range("I73").Value = 2089.3333
$ws.Range("K73").Value = 6267.999899999999
$ws.Range("M73").Value = -5175.999899999999
# Row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 1562.6875
$ws.Range("I107").Value = 1927.375
$ws.Range("J107").Value = 1198
$ws.Range("K107").Value = 5782.125
$ws.Range("L107").Value = 3594
$ws.Range("M107").Value = -3862.125
$ws.Range("N107").Value = -7434
# Row 113 (Leve Item ID 27843)
$ws.Range("H113").Value = 1139.24
$ws.Range("I113").Value = 686
$ws.Range("J113").Value = 1282.3684
$ws.Range("K113").Value = 2058
$ws.Range("L113").Value = 3847.1052
$ws.Range("M113").Value = 112
$ws.Range("N113").Value = -8187.1052
# Row 114 (Leve Item ID 27865)
$ws.Range("H114").Value = 691.9231
$ws.Range("I114").Value = 500
$ws.Range("J114").Value = 999
$ws.Range("K114").Value = 1500
$ws.Range("L114").Value = 2997
$ws.Range("M114").Value = 1754
$ws.Range("N114").Value = -9505
# Row 120 (Leve Item ID 27877)
$ws.Range("H120").Value = 5030
$ws.Range("I120").Value = 5030
$ws.Range("K120").Value = 15090
$ws.Range("M120").Value = -10252

$ws = $wb.Worksheets.Item("GSM")
# Row 43 (Leve Item ID 4218)
$ws.Range("H43").Value = 12224.889
$ws.Range("J43").Value = 20601.4
$ws.Range("L43").Value = 20601.4
$ws.Range("N43").Value = -20903.4
# Row 44 (Leve Item ID 4143)
$ws.Range("H44").Value = 4312476
$ws.Range("I44").Value = 33166.5
$ws.Range("K44").Value = 33166.5
$ws.Range("M44").Value = -32570.5
# Row 92 (Leve Item ID 18094)
$ws.Range("H92").Value = 21500
$ws.Range("J92").Value = 21500
$ws.Range("L92").Value = 21500
$ws.Range("N92").Value = -25244
# Row 97 (Leve Item ID 19940)
$ws.Range("H97").Value = 1256.5897
$ws.Range("J97").Value = 4897.5
$ws.Range("L97").Value = 4897.5
$ws.Range("N97").Value = -5889.5
# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 2982.1667
$ws.Range("I126").Value = 2793.5557
$ws.Range("K126").Value = 8380.667099999999
$ws.Range("M126").Value = -5910.667099999999
# Row 140 (Leve Item ID 42458)
$ws.Range("H140").Value = 189750
$ws.Range("J140").Value = 189750
$ws.Range("L140").Value = 189750
$ws.Range("N140").Value = -200110

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 3721
$ws.Range("I7").Value = 3712.7144
$ws.Range("J7").Value = 3750
$ws.Range("K7").Value = 3712.7144
$ws.Range("L7").Value = 3750
$ws.Range("M7").Value = -3600.7144
$ws.Range("N7").Value = -3974
# Row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 2065.6667
$ws.Range("I16").Value = 2252.125
$ws.Range("J16").Value = 1692.75
$ws.Range("K16").Value = 2252.125
$ws.Range("L16").Value = 1692.75
$ws.Range("M16").Value = -2082.125
$ws.Range("N16").Value = -2032.75
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 3991.875
$ws.Range("I40").Value = 3992.1428
$ws.Range("K40").Value = 3992.1428
$ws.Range("M40").Value = -3856.1428
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 4801.5
$ws.Range("I46").Value = 1485
$ws.Range("K46").Value = 1485
$ws.Range("M46").Value = -1297
# Row 82 (Leve Item ID 12565)
$ws.Range("H82").Value = 11027.48
$ws.Range("I82").Value = 10271.866
$ws.Range("K82").Value = 10271.866
$ws.Range("M82").Value = -9910.866
# Row 85 (Leve Item ID 12565)
$ws.Range("H85").Value = 11027.48
$ws.Range("I85").Value = 10271.866
$ws.Range("K85").Value = 10271.866
$ws.Range("M85").Value = -9023.866
# Row 100 (Leve Item ID 19995)
$ws.Range("H100").Value = 2934.95
$ws.Range("I100").Value = 2811.7058
$ws.Range("K100").Value = 2811.7058
$ws.Range("M100").Value = -2270.7058
# Row 108 (Leve Item ID 25655)
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 3721
$ws.Range("I126").Value = 3712.7144
$ws.Range("J126").Value = 3750
$ws.Range("K126").Value = 11138.1432
$ws.Range("L126").Value = 11250
$ws.Range("M126").Value = -8668.143199999999
$ws.Range("N126").Value = -16190

$ws = $wb.Worksheets.Item("WVR")
# Row 141 (Leve Item ID 42505)
$ws.Range("H141").Value = 102593.2
$ws.Range("J141").Value = 95822.336
$ws.Range("L141").Value = 95822.336
$ws.Range("N141").Value = -106182.336

Write-Host "Applied scheduled market refresh updates"